{"js": "// Remove the \"$> pip install faiss-gpu==1.7.2\" command-line paragraph that\n// immediately follows the \"$> pip install onnxruntime-gpu==1.16.3\" paragraph\n// in the conflicting-dependencies install steps. (There is a second,\n// unrelated \"pip install faiss-gpu==1.7.2\" paragraph later in the document,\n// under the \"Install the GPU version of faiss\" heading, which must be left\n// untouched.)\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst NEEDLE = \"pip install faiss-gpu==1.7.2\";\nconst ANCHOR = \"pip install onnxruntime-gpu==1.16.3\";\n\nlet target = null;\nfor (let i = 1; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.indexOf(NEEDLE) !== -1 &&\n      paragraphs.items[i - 1].text.indexOf(ANCHOR) !== -1) {\n    target = para;\n    break;\n  }\n}\n\nif (target) {\n  target.delete();\n  await context.sync();\n}\n", "ps1": "# Remove the \"$> pip install faiss-gpu==1.7.2\" command-line paragraph that\n# immediately follows the \"$> pip install onnxruntime-gpu==1.16.3\" paragraph\n# in the conflicting-dependencies install steps. (There is a second,\n# unrelated \"pip install faiss-gpu==1.7.2\" paragraph later in the document,\n# under the \"Install the GPU version of faiss\" heading, which must be left\n# untouched.)\n\n$d = $word.ActiveDocument\n\n$wdParagraph = 4\n$wdCollapseEnd = 0\n\n# Locate the unique anchor line.\n$range = $d.Content\n$find = $range.Find\n$find.Text = \"pip install onnxruntime-gpu==1.16.3\"\n$found = $find.Execute()\n\nif ($found) {\n    # Grow the found hit to its whole paragraph (including the paragraph\n    # mark), collapse to the end of it, then grow again to select the\n    # following paragraph (including its paragraph mark).\n    $range.Expand($wdParagraph) | Out-Null\n    $range.Collapse($wdCollapseEnd)\n    $range.Expand($wdParagraph) | Out-Null\n\n    if ($range.Text -like \"*pip install faiss-gpu==1.7.2*\") {\n        $range.Delete()\n    }\n}\n"}
